$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 formatting: remove the yellow highlight that was applied to
# "Servo Claw" / "Drop Mechanism2" so it matches the rest of the table
# (copy the plain formatting used by the other data rows, e.g. row 12,
# onto row 4 without touching its values/formula).
$ws.Range("A12:G12").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Update the unit price for the "Servo Claw" line; the dependent
# Total formula (G4) and the subtotal (G38) recalc automatically.
$ws.Range("F4").Value = 22.9

# --- Sheet view: drop the frozen/top-left scroll position and move the
# active selection to C9.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C9").Select()
